$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, centered, bordered) from existing header cell (H1) to new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for new columns I (I0) and J (IF)
$iValues = @{2=1; 3=1; 4=1; 5=1; 6=1; 7=1; 8=1; 9=1; 10=1; 11=4; 12=1}
$jValues = @{2=4; 3=7; 4=5; 5=6; 6=5; 7=4; 8=5; 9=4; 10=3; 11=5; 12=1}

foreach ($r in 2..12) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
